# chore: update Sheets via scheduled runner
# Refreshes market-price derived figures (currentAveragePrice*, LevePrice*,
# LeveProfit*) across the per-job leve-profit worksheets (ALC, ARM, BSM, CRP,
# CUL, GSM, LTW, WVR) to reflect newly pulled market data.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 166.8
$ws.Range("I33").Value = 88.75
$ws.Range("K33").Value = 88.75
$ws.Range("M33").Value = 140.25
$ws.Range("H76").Value = 6198.25
$ws.Range("I76").Value = 6076.5
$ws.Range("K76").Value = 6076.5
$ws.Range("M76").Value = -5761.5
$ws.Range("H79").Value = 6198.25
$ws.Range("I79").Value = 6076.5
$ws.Range("K79").Value = 6076.5
$ws.Range("M79").Value = -4984.5
$ws.Range("H106").Value = 13533.8
$ws.Range("I106").Value = 15854
$ws.Range("K106").Value = 15854
$ws.Range("M106").Value = -15223
$ws.Range("H107").Value = 1905.2106
$ws.Range("I107").Value = 1706.4615
$ws.Range("J107").Value = 2335.8333
$ws.Range("K107").Value = 1706.4615
$ws.Range("L107").Value = 2335.8333
$ws.Range("M107").Value = 213.5385000000001
$ws.Range("N107").Value = -6175.8333
$ws.Range("H132").Value = 11501698
$ws.Range("I132").Value = 15157538
$ws.Range("K132").Value = 45472614
$ws.Range("M132").Value = -45470084
$ws.Range("H137").Value = 1289.7037
$ws.Range("J137").Value = 2065.7778
$ws.Range("L137").Value = 6197.3334
$ws.Range("N137").Value = -11297.3334
$ws.Range("H138").Value = 1577.5745
$ws.Range("J138").Value = 1922.4062
$ws.Range("L138").Value = 5767.2186
$ws.Range("N138").Value = -16047.2186
$ws.Range("H141").Value = 400
$ws.Range("I141").Value = 400
$ws.Range("J141").Value = 0
$ws.Range("K141").Value = 1200
$ws.Range("L141").Value = 0
$ws.Range("M141").Value = 3980
$ws.Range("N141").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 9933.091
$ws.Range("I2").Value = 560.1429000000001
$ws.Range("J2").Value = 26335.75
$ws.Range("K2").Value = 560.1429000000001
$ws.Range("L2").Value = 26335.75
$ws.Range("M2").Value = -447.1429000000001
$ws.Range("N2").Value = -26561.75
$ws.Range("H32").Value = 2948.3376
$ws.Range("I32").Value = 2857.7014
$ws.Range("J32").Value = 3555.6
$ws.Range("K32").Value = 2857.7014
$ws.Range("L32").Value = 3555.6
$ws.Range("M32").Value = -2570.7014
$ws.Range("N32").Value = -4129.6
$ws.Range("H45").Value = 1112.6
$ws.Range("I45").Value = 1168.9
$ws.Range("J45").Value = 1000
$ws.Range("K45").Value = 1168.9
$ws.Range("L45").Value = 1000
$ws.Range("M45").Value = -791.9000000000001
$ws.Range("N45").Value = -1754
$ws.Range("H61").Value = 1307.75
$ws.Range("I61").Value = 1077
$ws.Range("J61").Value = 2000
$ws.Range("K61").Value = 1077
$ws.Range("L61").Value = 2000
$ws.Range("M61").Value = -865
$ws.Range("N61").Value = -2424
$ws.Range("H74").Value = 813.3570999999999
$ws.Range("I74").Value = 645.1539
$ws.Range("J74").Value = 3000
$ws.Range("K74").Value = 645.1539
$ws.Range("L74").Value = 3000
$ws.Range("M74").Value = 228.8461
$ws.Range("N74").Value = -4748
$ws.Range("H77").Value = 813.3570999999999
$ws.Range("I77").Value = 645.1539
$ws.Range("J77").Value = 3000
$ws.Range("K77").Value = 3225.7695
$ws.Range("L77").Value = 15000
$ws.Range("M77").Value = 1142.2305
$ws.Range("N77").Value = -23736
$ws.Range("H96").Value = 0
$ws.Range("J96").Value = 0
$ws.Range("L96").Value = 0
$ws.Range("N96").ClearContents()
$ws.Range("H110").Value = 1093.3334
$ws.Range("I110").Value = 930
$ws.Range("J110").Value = 1501.6666
$ws.Range("K110").Value = 930
$ws.Range("L110").Value = 1501.6666
$ws.Range("M110").Value = 1115
$ws.Range("N110").Value = -5591.6666
$ws.Range("H116").Value = 9933.091
$ws.Range("I116").Value = 560.1429000000001
$ws.Range("J116").Value = 26335.75
$ws.Range("K116").Value = 560.1429000000001
$ws.Range("L116").Value = 26335.75
$ws.Range("M116").Value = 1733.8571
$ws.Range("N116").Value = -30923.75
$ws.Range("H136").Value = 1307.75
$ws.Range("I136").Value = 1077
$ws.Range("J136").Value = 2000
$ws.Range("K136").Value = 3231
$ws.Range("L136").Value = 6000
$ws.Range("M136").Value = -681
$ws.Range("N136").Value = -11100

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 9933.091
$ws.Range("I3").Value = 560.1429000000001
$ws.Range("J3").Value = 26335.75
$ws.Range("K3").Value = 560.1429000000001
$ws.Range("L3").Value = 26335.75
$ws.Range("M3").Value = -446.1429000000001
$ws.Range("N3").Value = -26563.75
$ws.Range("H20").Value = 1517
$ws.Range("I20").Value = 1104.75
$ws.Range("K20").Value = 1104.75
$ws.Range("M20").Value = -857.75
$ws.Range("H134").Value = 7749.278
$ws.Range("I134").Value = 1204.9166
$ws.Range("J134").Value = 20838
$ws.Range("K134").Value = 3614.7498
$ws.Range("L134").Value = 62514
$ws.Range("M134").Value = -1079.7498
$ws.Range("N134").Value = -67584

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 52632504
$ws.Range("I16").Value = 66667620
$ws.Range("J16").Value = 814.25
$ws.Range("K16").Value = 66667620
$ws.Range("L16").Value = 814.25
$ws.Range("M16").Value = -66667333
$ws.Range("N16").Value = -1388.25
$ws.Range("H31").Value = 1354.9803
$ws.Range("I31").Value = 1308.2449
$ws.Range("J31").Value = 2500
$ws.Range("K31").Value = 1308.2449
$ws.Range("L31").Value = 2500
$ws.Range("M31").Value = -1013.2449
$ws.Range("N31").Value = -3090
$ws.Range("H34").Value = 1354.9803
$ws.Range("I34").Value = 1308.2449
$ws.Range("J34").Value = 2500
$ws.Range("K34").Value = 1308.2449
$ws.Range("L34").Value = 2500
$ws.Range("M34").Value = -1106.2449
$ws.Range("N34").Value = -2904
$ws.Range("H58").Value = 1531
$ws.Range("I58").Value = 1339.7778
$ws.Range("K58").Value = 1339.7778
$ws.Range("M58").Value = -1136.7778
$ws.Range("H113").Value = 52632504
$ws.Range("I113").Value = 66667620
$ws.Range("J113").Value = 814.25
$ws.Range("K113").Value = 66667620
$ws.Range("L113").Value = 814.25
$ws.Range("M113").Value = -66665450
$ws.Range("N113").Value = -5154.25
$ws.Range("H122").Value = 758.4737
$ws.Range("I122").Value = 687.4
$ws.Range("K122").Value = 2062.2
$ws.Range("M122").Value = 387.8000000000002
$ws.Range("H132").Value = 2761.889
$ws.Range("I132").Value = 2065.0908
$ws.Range("J132").Value = 3856.8572
$ws.Range("K132").Value = 6195.2724
$ws.Range("L132").Value = 11570.5716
$ws.Range("M132").Value = -3665.2724
$ws.Range("N132").Value = -16630.5716
$ws.Range("H134").Value = 1319
$ws.Range("I134").Value = 1188.909
$ws.Range("K134").Value = 3566.727
$ws.Range("M134").Value = -1031.727
$ws.Range("H136").Value = 1531
$ws.Range("I136").Value = 1339.7778
$ws.Range("K136").Value = 4019.3334
$ws.Range("M136").Value = -1469.3334

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H87").Value = 2951.375
$ws.Range("I87").Value = 1604
$ws.Range("J87").Value = 3759.8
$ws.Range("K87").Value = 4812
$ws.Range("L87").Value = 11279.4
$ws.Range("M87").Value = -3564
$ws.Range("N87").Value = -13775.4
$ws.Range("H90").Value = 2951.375
$ws.Range("I90").Value = 1604
$ws.Range("J90").Value = 3759.8
$ws.Range("K90").Value = 14436
$ws.Range("L90").Value = 33838.2
$ws.Range("M90").Value = -8196
$ws.Range("N90").Value = -46318.2
$ws.Range("H131").Value = 13335621
$ws.Range("J131").Value = 2655.746
$ws.Range("L131").Value = 7967.238
$ws.Range("N131").Value = -18047.238

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 6078
$ws.Range("I102").Value = 6078
$ws.Range("K102").Value = 6078
$ws.Range("M102").Value = -4456

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 761.0714
$ws.Range("J22").Value = 878.44446
$ws.Range("L22").Value = 878.44446
$ws.Range("N22").Value = -1468.44446
$ws.Range("H27").Value = 761.0714
$ws.Range("J27").Value = 878.44446
$ws.Range("L27").Value = 878.44446
$ws.Range("N27").Value = -1092.44446
$ws.Range("H40").Value = 3088.7856
$ws.Range("I40").Value = 2879
$ws.Range("K40").Value = 2879
$ws.Range("M40").Value = -2743
$ws.Range("H46").Value = 2750
$ws.Range("I46").Value = 1000
$ws.Range("J46").Value = 3333.3333
$ws.Range("K46").Value = 1000
$ws.Range("L46").Value = 3333.3333
$ws.Range("M46").Value = -812
$ws.Range("N46").Value = -3709.3333
$ws.Range("H122").Value = 35716470
$ws.Range("I122").Value = 50002176
$ws.Range("J122").Value = 2205
$ws.Range("K122").Value = 150006528
$ws.Range("L122").Value = 6615
$ws.Range("M122").Value = -150004078
$ws.Range("N122").Value = -11515
$ws.Range("H139").Value = 52810
$ws.Range("J139").Value = 52810
$ws.Range("L139").Value = 52810
$ws.Range("N139").Value = -63090

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H31").Value = 29000
$ws.Range("J31").Value = 29000
$ws.Range("L31").Value = 29000
$ws.Range("N31").Value = -29696
$ws.Range("H107").Value = 409.1579
$ws.Range("I107").Value = 358.26666
$ws.Range("J107").Value = 600
$ws.Range("K107").Value = 1074.79998
$ws.Range("L107").Value = 1800
$ws.Range("M107").Value = 845.20002
$ws.Range("N107").Value = -5640
$ws.Range("H117").Value = 25000
$ws.Range("J117").Value = 25000
$ws.Range("L117").Value = 25000
$ws.Range("N117").Value = -34178
$ws.Range("H136").Value = 493.63635
$ws.Range("I136").Value = 493.63635
$ws.Range("K136").Value = 1480.90905
$ws.Range("M136").Value = 1069.09095
